# Update cosinor per-day fixed-period-0 stats after rerunning the
# CircaDB / CircadiPy analyses (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[53.139511656388734, 71.73797017985035]"
$ws.Range("U2").Value = "[42.709239773808704, 55.05502837192985]"
$ws.Range("M3").Value = "[50.474132317907134, 75.07620669564143]"
$ws.Range("N3").Value = [double]"2.184918912462308e-13"
$ws.Range("O3").Value = [double]"2.184918912462308e-13"
$ws.Range("Q3").Value = "[0.9497106920761942, 1.3522370781217328]"
$ws.Range("R3").Value = [double]"5.10702591327572e-15"
$ws.Range("S3").Value = [double]"5.10702591327572e-15"
$ws.Range("U3").Value = "[42.09377364043489, 55.552760463525246]"
$ws.Range("Y3").Value = [double]"20.23959959960007"
$ws.Range("Z3").Value = [double]"21.89181181181231"
$ws.Range("M4").Value = "[51.07984507923581, 78.04827029292073]"
$ws.Range("N4").Value = [double]"1.599831378484851e-12"
$ws.Range("O4").Value = [double]"1.599831378484851e-12"
$ws.Range("U4").Value = "[43.30976594730957, 56.831947760793966]"
$ws.Range("M5").Value = "[54.0400621764807, 77.80534601532689]"
$ws.Range("N5").Value = [double]"1.443289932012704e-14"
$ws.Range("O5").Value = [double]"1.443289932012704e-14"
$ws.Range("Q5").Value = "[0.10692107129334527, 0.45913165908319353]"
$ws.Range("R5").Value = [double]"0.002270104082773727"
$ws.Range("S5").Value = [double]"0.002270104082773727"
$ws.Range("U5").Value = "[44.022516316061015, 57.40579830892854]"
$ws.Range("Y5").Value = [double]"23.90544544544599"
$ws.Range("Z5").Value = [double]"25.35113113113172"
$ws.Range("M6").Value = "[54.595203645818614, 75.63858549460672]"
$ws.Range("N6").Value = [double]"2.220446049250313e-16"
$ws.Range("O6").Value = [double]"2.220446049250313e-16"
$ws.Range("Q6").Value = "[-0.2515789912784623, 0.10063159651138598]"
$ws.Range("R6").Value = [double]"0.3926119510892936"
$ws.Range("S6").Value = [double]"0.3926119510892936"
$ws.Range("U6").Value = "[43.60683008097774, 57.091174378075245]"
$ws.Range("Y6").Value = [double]"-0.4130530530530662"
$ws.Range("Z6").Value = [double]"1.032632632632659"
$ws.Range("M7").Value = "[49.91078193105742, 75.63970082452781]"
$ws.Range("N7").Value = [double]"8.921752225887758e-13"
$ws.Range("O7").Value = [double]"8.921752225887758e-13"
$ws.Range("U7").Value = "[42.09294946522694, 55.55321874949185]"
$ws.Range("M8").Value = "[50.58931241449935, 73.95794283260423]"
$ws.Range("N8").Value = [double]"5.395683899678261e-14"
$ws.Range("O8").Value = [double]"5.395683899678261e-14"
$ws.Range("U8").Value = "[40.869400269977845, 53.320543734657704]"
$ws.Range("M9").Value = "[53.9326023516647, 73.25287361249023]"
$ws.Range("U9").Value = "[43.83696765073083, 56.163270446630705]"
$ws.Range("M10").Value = "[53.88980850828719, 72.73071627766234]"
$ws.Range("U10").Value = "[43.44748942536745, 56.07627568482188]"
$ws.Range("M11").Value = "[53.578942350451584, 72.89390156202124]"
$ws.Range("Q11").Value = "[-1.7987897876410024, -1.446579199851156]"
$ws.Range("U11").Value = "[45.39995869560904, 58.062498928216094]"
$ws.Range("Y11").Value = [double]"5.679779779779875"
$ws.Range("Z11").Value = [double]"7.062682682682801"
